$d = $word.ActiveDocument

# 1. Replace the ">>> your stuff after this line >>>" paragraph text
#    (3 runs) with a single run "A quote by Eknath Daster:"
$d.Content.Find.Execute(">>> your stuff after this line >>>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A quote by Eknath Daster:", 2)

# 2. Insert a new paragraph right after the "A quote by Eknath Daster:" paragraph (P4)
#    for "Programming is thinking, not typing."
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()

$typingPara = $d.Paragraphs(5)
# type the text plus a trailing sentinel character so the later bookmark
# insertion point is never the very-last-character-of-paragraph edge case
$typingPara.Range.Text = "Programming is thinking, not typing.#"

$typingPara2 = $d.Paragraphs(5)
$bnd = $typingPara2.Range.End - 2   # just before the sentinel '#', i.e. right after "."
$bmRange = $d.Range($bnd, $bnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# remove the sentinel character now that the bookmark is anchored
$typingPara3 = $d.Paragraphs(5)
$sentinelRange = $d.Range($typingPara3.Range.End - 2, $typingPara3.Range.End - 1)
$sentinelRange.Delete()

# 3. Insert the blank paragraph before the "A quote by Eknath Daster:" paragraph
#    (i.e. right after paragraph 3, "...assignment report!")
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()

# 4. Insert the blank paragraph after the typing paragraph (i.e. right before
#    the "A cool quote by Dijkstra:" paragraph). After step 3 all subsequent
#    paragraph indices shifted by +1, so the typing paragraph is now #6.
$typingParaFinal = $d.Paragraphs(6)
$typingParaFinal.Range.InsertParagraphAfter()
